$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "20.243.91"
Set-TextValue $ws.Range("E2") "  +1.83%  "
Set-TextValue $ws.Range("D3") "1.441.35"
Set-TextValue $ws.Range("E3") "  +3.09%  "
Set-TextValue $ws.Range("D4") "1.006"
Set-TextValue $ws.Range("E4") "  +0.22%  "
Set-TextValue $ws.Range("D5") "0.9157"
Set-TextValue $ws.Range("E5") "  -8.80%  "
Set-TextValue $ws.Range("D6") "274.71"
Set-TextValue $ws.Range("E6") "  +0.58%  "
Set-TextValue $ws.Range("E7") "  -1.10%  "
Set-TextValue $ws.Range("D8") "0.3077"
Set-TextValue $ws.Range("E8") "  -0.91%  "
Set-TextValue $ws.Range("D9") "39.04"
Set-TextValue $ws.Range("E9") "  -1.27%  "
Set-TextValue $ws.Range("D10") "1.019"
Set-TextValue $ws.Range("E10") "  +1.79%  "
Set-TextValue $ws.Range("D11") "0.06497"
Set-TextValue $ws.Range("E11") "  +0.46%  "
Set-TextValue $ws.Range("D12") "0.9993"
Set-TextValue $ws.Range("E12") "  -0.46%  "
Set-TextValue $ws.Range("D13") "5.338"
Set-TextValue $ws.Range("E13") "  -1.40%  "
Set-TextValue $ws.Range("D14") "17.45"
Set-TextValue $ws.Range("E14") "  +1.30%  "
Set-TextValue $ws.Range("D15") "6.034"
Set-TextValue $ws.Range("E15") "  -1.38%  "
Set-TextValue $ws.Range("D16") "0.00001009"
Set-TextValue $ws.Range("E16") "  +0.01%  "
Set-TextValue $ws.Range("D17") "1.441.66"
Set-TextValue $ws.Range("E17") "  +3.12%  "
Set-TextValue $ws.Range("D18") "0.9310"
Set-TextValue $ws.Range("E18") "  -7.29%  "
Set-TextValue $ws.Range("E19") "  -1.08%  "
Set-TextValue $ws.Range("D20") "67.42"
Set-TextValue $ws.Range("E20") "  -3.57%  "
Set-TextValue $ws.Range("D21") "5.401"
Set-TextValue $ws.Range("E21") "  -2.58%  "
Set-TextValue $ws.Range("D22") "14.19"
Set-TextValue $ws.Range("E22") "  -2.81%  "
Set-TextValue $ws.Range("D23") "10.81"
Set-TextValue $ws.Range("E23") "  -1.52%  "
Set-TextValue $ws.Range("D24") "2.229"
Set-TextValue $ws.Range("E24") "  -1.73%  "
Set-TextValue $ws.Range("D25") "20.280.58"
Set-TextValue $ws.Range("E25") "  +2.04%  "
Set-TextValue $ws.Range("D26") "137.66"
Set-TextValue $ws.Range("E26") "  +1.99%  "
Set-TextValue $ws.Range("D27") "2.059"
Set-TextValue $ws.Range("E27") "  -6.68%  "
Set-TextValue $ws.Range("D28") "16.90"
Set-TextValue $ws.Range("E28") "  +0.40%  "
Set-TextValue $ws.Range("D29") "1.593.38"
Set-TextValue $ws.Range("E29") "  +2.31%  "
Set-TextValue $ws.Range("D30") "110.06"
Set-TextValue $ws.Range("E30") "  +1.06%  "
Set-TextValue $ws.Range("D31") "3.979"
Set-TextValue $ws.Range("E31") "  -2.22%  "
Set-TextValue $ws.Range("D32") "0.7966"
Set-TextValue $ws.Range("E32") "  -1.07%  "
Set-TextValue $ws.Range("D33") "4.831"
Set-TextValue $ws.Range("E33") "  -7.97%  "
Set-TextValue $ws.Range("D34") "0.07661"
Set-TextValue $ws.Range("E34") "  +0.46%  "
Set-TextValue $ws.Range("D35") "1.464"
Set-TextValue $ws.Range("E35") "  +2.57%  "
Set-TextValue $ws.Range("D36") "0.05782"
Set-TextValue $ws.Range("E36") "  +0.49%  "
Set-TextValue $ws.Range("D37") "4.657"
Set-TextValue $ws.Range("E37") "  -2.72%  "
Set-TextValue $ws.Range("D38") "1.136"
Set-TextValue $ws.Range("E38") "  +4.15%  "
Set-TextValue $ws.Range("D39") "0.01983"
Set-TextValue $ws.Range("E39") "  -3.45%  "
Set-TextValue $ws.Range("D40") "10.13"
Set-TextValue $ws.Range("E40") "  -1.40%  "
Set-TextValue $ws.Range("D41") "0.1848"
Set-TextValue $ws.Range("E41") "  -2.00%  "
Set-TextValue $ws.Range("D42") "0.9272"
Set-TextValue $ws.Range("E42") "  -7.58%  "
Set-TextValue $ws.Range("D43") "7.014"
Set-TextValue $ws.Range("E43") "  -15.88%  "
Set-TextValue $ws.Range("D44") "0.5192"
Set-TextValue $ws.Range("E44") "  -0.80%  "
Set-TextValue $ws.Range("D45") "3.479"
Set-TextValue $ws.Range("E45") "  -0.80%  "
Set-TextValue $ws.Range("D46") "11.78"
Set-TextValue $ws.Range("E46") "  -3.23%  "
Set-TextValue $ws.Range("D47") "116.24"
Set-TextValue $ws.Range("E47") "  +4.63%  "
Set-TextValue $ws.Range("D48") "0.5102"
Set-TextValue $ws.Range("E48") "  +0.76%  "
Set-TextValue $ws.Range("D49") "1.724"
Set-TextValue $ws.Range("E49") "  -1.27%  "
Set-TextValue $ws.Range("D50") "0.06406"
Set-TextValue $ws.Range("E50") "  +4.48%  "
Set-TextValue $ws.Range("D51") "0.9740"
Set-TextValue $ws.Range("E51") "  -3.26%  "
